$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark LED row (row 6) as Do Not Install, and swap the vendor/mfg part to the
# Lite-On replacement (was Cree).
$ws.Range("C6").Value = "DNI"
$ws.Range("G6").Value = "160-1940-ND"
$ws.Range("H6").Value = "LTL-1CHA"
$ws.Range("I6").Value = "Lite-On"

# Mark the rotary encoder row (row 22) as Do Not Install too.
$ws.Range("C22").Value = "DNI"

# Move the active selection to B6.
$ws.Range("B6").Select() | Out-Null
